$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 8.3012195145644281
$ws.Range("C2").Value = 5.5808170512609685
$ws.Range("D2").Value = 9.013385994064695
$ws.Range("E2").Value = 7.4874158363080605

$ws.Range("B3").Value = 5.7196125950586341
$ws.Range("C3").Value = 7.025589783105973
$ws.Range("D3").Value = 5.3717862007717372
$ws.Range("E3").Value = 8.3214925567711493

$ws.Range("B1:E3").Select() | Out-Null
